# Apply the update:
#  - Column C (row 2..16): 46062 -> 46063 ("Förändrad" date +1 day)
#  - Rows 8, 10, 12, 14, 15: re-ordered records (A/B/G values moved between rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the "Förändrad" (changed) date column C for every data row.
for ($r = 2; $r -le 16; $r++) {
    $ws.Range("C$r").Value = 46063
}

# Row re-ordering among rows 8, 10, 12, 14, 15 (A, B, G columns only).
$ws.Range("A8").Value  = "A 50530-2024"
$ws.Range("B8").Value  = 45601.56424768519
$ws.Range("G8").Value  = 0.7

$ws.Range("A10").Value = "A 23678-2023"
$ws.Range("B10").Value = 45077
$ws.Range("G10").Value = 1.4

$ws.Range("A12").Value = "A 50538-2024"
$ws.Range("B12").Value = 45601.57153935185
$ws.Range("G12").Value = 0.8

$ws.Range("A14").Value = "A 45370-2022"
$ws.Range("B14").Value = 44844.6397337963
$ws.Range("G14").Value = 2.7

$ws.Range("A15").Value = "A 2253-2022"
$ws.Range("B15").Value = 44578
$ws.Range("G15").Value = 0.3
